$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (F1:H1) - styled like the existing headers (bold, centered, bordered)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the style of an existing header cell (A1) onto the new header cells
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Boolean outlier flags for rows 2-9
$values = @(
    @($false, $false, $false),  # row 2
    @($false, $false, $false),  # row 3
    @($false, $false, $false),  # row 4
    @($true,  $false, $false),  # row 5
    @($false, $false, $false),  # row 6
    @($false, $false, $false),  # row 7
    @($false, $false, $false),  # row 8
    @($true,  $false, $false)   # row 9
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i][0]
    $ws.Cells.Item($row, 7).Value = $values[$i][1]
    $ws.Cells.Item($row, 8).Value = $values[$i][2]
}
